$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Cells.Item(2, 3).Value = 0.0259985220088339
$ws.Cells.Item(2, 4).Value = 0.02370252457615862
$ws.Cells.Item(2, 5).Value = 2.610317677428952
$ws.Cells.Item(2, 6).Value = 0.2073471222403143
$ws.Cells.Item(2, 7).Value = 0.09184754780432058
$ws.Cells.Item(2, 8).Value = 0.2387076879951664
$ws.Cells.Item(2, 9).Value = 0.1444412836853388
$ws.Cells.Item(2, 13).Value = 10.73723423822503
$ws.Cells.Item(2, 15).Value = 0.5480276828193098
$ws.Cells.Item(3, 3).Value = 0.02290818794315896
$ws.Cells.Item(3, 4).Value = 0.02072920476602746
$ws.Cells.Item(3, 5).Value = 2.273873653202685
$ws.Cells.Item(3, 6).Value = 0.2094294264737115
$ws.Cells.Item(3, 7).Value = 0.09279224558296306
$ws.Cells.Item(3, 8).Value = 0.2472243309196216
$ws.Cells.Item(3, 9).Value = 0.1467012107897148
$ws.Cells.Item(3, 13).Value = 9.377885926079102
$ws.Cells.Item(3, 15).Value = 0.5674014379640084
$ws.Cells.Item(4, 3).Value = 0.02100197363344591
$ws.Cells.Item(4, 4).Value = 0.01889398304664525
$ws.Cells.Item(4, 5).Value = 2.067786694163914
$ws.Cells.Item(4, 6).Value = 0.2113172490311257
$ws.Cells.Item(4, 7).Value = 0.09390085873008758
$ws.Cells.Item(4, 8).Value = 0.2529030733522788
$ws.Cells.Item(4, 9).Value = 0.1485473587455495
$ws.Cells.Item(4, 13).Value = 8.541643043153556
$ws.Cells.Item(4, 15).Value = 0.5813463997372708
$ws.Cells.Item(5, 3).Value = 0.02022307162040704
$ws.Cells.Item(5, 4).Value = 0.018143751160558
$ws.Cells.Item(5, 5).Value = 1.983901584865976
$ws.Cells.Item(5, 6).Value = 0.2122357318661727
$ws.Cells.Item(5, 7).Value = 0.09448073276683999
$ws.Cells.Item(5, 8).Value = 0.2553282181036067
$ws.Cells.Item(5, 9).Value = 0.1494118680006977
$ws.Cells.Item(5, 13).Value = 8.200391630924969
$ws.Cells.Item(5, 15).Value = 0.5875299969409866
$ws.Cells.Item(6, 3).Value = 0.02009361126150111
$ws.Cells.Item(6, 4).Value = 0.01801903452156495
$ws.Cells.Item(6, 5).Value = 1.96997765980305
$ws.Cells.Item(6, 6).Value = 0.2123971385254322
$ws.Cells.Item(6, 7).Value = 0.09458461854526945
$ws.Cells.Item(6, 8).Value = 0.2557375577692937
$ws.Cells.Item(6, 9).Value = 0.1495621047293376
$ws.Cells.Item(6, 13).Value = 8.143696359199225
$ws.Cells.Item(6, 15).Value = 0.5885866145087135
$ws.Cells.Item(7, 3).Value = 0.02099147748216978
$ws.Cells.Item(7, 4).Value = 0.01888387465977104
$ws.Cells.Item(7, 5).Value = 2.066655030645848
$ws.Cells.Item(7, 6).Value = 0.2113290372707155
$ws.Cells.Item(7, 7).Value = 0.09390816671511004
$ws.Cells.Item(7, 8).Value = 0.2529353328685744
$ws.Cells.Item(7, 9).Value = 0.1485585676551118
$ws.Cells.Item(7, 13).Value = 8.537042825085166
$ws.Cells.Item(7, 15).Value = 0.5814277846186258
$ws.Cells.Item(8, 3).Value = 0.02493484135088408
$ws.Cells.Item(8, 4).Value = 0.02267933938777844
$ws.Cells.Item(8, 5).Value = 2.494193039416473
$ws.Cells.Item(8, 6).Value = 0.2079362038840316
$ws.Cells.Item(8, 7).Value = 0.09206065931312679
$ws.Cells.Item(8, 8).Value = 0.2415497160192501
$ws.Cells.Item(8, 9).Value = 0.1451234478444832
$ws.Cells.Item(8, 13).Value = 10.26881511767897
$ws.Cells.Item(8, 15).Value = 0.5542736574510059
$ws.Cells.Item(9, 3).Value = 0.032594964002854
$ws.Cells.Item(9, 4).Value = 0.03004489036957381
$ws.Cells.Item(9, 5).Value = 3.337856558203612
$ws.Cells.Item(9, 6).Value = 0.2063010013012345
$ws.Cells.Item(9, 7).Value = 0.0928514382002561
$ws.Cells.Item(9, 8).Value = 0.2228802192032333
$ws.Cells.Item(9, 9).Value = 0.1421676255902149
$ws.Cells.Item(9, 13).Value = 13.65605169794696
$ws.Cells.Item(9, 15).Value = 0.5179430773937526
$ws.Cells.Item(10, 3).Value = 0.03817466801056923
$ws.Cells.Item(10, 4).Value = 0.03540883395066885
$ws.Cells.Item(10, 5).Value = 3.963061294130966
$ws.Cells.Item(10, 6).Value = 0.2084172051241424
$ws.Cells.Item(10, 7).Value = 0.09643301549894545
$ws.Cells.Item(10, 8).Value = 0.2115205661281365
$ws.Cells.Item(10, 9).Value = 0.1425004041224511
$ws.Cells.Item(10, 13).Value = 16.14560305806549
$ws.Cells.Item(10, 15).Value = 0.5024888564499435
$ws.Cells.Item(11, 3).Value = 0.04070188383830953
$ws.Cells.Item(11, 4).Value = 0.03783889709913524
$ws.Cells.Item(11, 5).Value = 4.249185403936622
$ws.Cells.Item(11, 6).Value = 0.2101579567766265
$ws.Cells.Item(11, 7).Value = 0.09878337973663065
$ws.Cells.Item(11, 8).Value = 0.206892742323646
$ws.Cells.Item(11, 9).Value = 0.143240205404858
$ws.Cells.Item(11, 13).Value = 17.27995216389007
$ws.Cells.Item(11, 15).Value = 0.4981050916419463
$ws.Cells.Item(12, 3).Value = 0.04165722912055969
$ws.Cells.Item(12, 4).Value = 0.0387576772922813
$ws.Cells.Item(12, 5).Value = 4.3578279018069
$ws.Cells.Item(12, 6).Value = 0.2109339253391909
$ws.Cells.Item(12, 7).Value = 0.09978303283034506
$ws.Cells.Item(12, 8).Value = 0.2052203371492354
$ws.Cells.Item(12, 9).Value = 0.1436087631029039
$ws.Cells.Item(12, 13).Value = 17.70990585949443
$ws.Cells.Item(12, 15).Value = 0.4968433022546037
$ws.Cells.Item(13, 3).Value = 0.04145155318222749
$ws.Cells.Item(13, 4).Value = 0.03855986494727404
$ws.Cells.Item(13, 5).Value = 4.334415951801134
$ws.Cells.Item(13, 6).Value = 0.210761533037207
$ws.Cells.Item(13, 7).Value = 0.09956276763634975
$ws.Cells.Item(13, 8).Value = 0.2055769196510084
$ws.Cells.Item(13, 9).Value = 0.1435253937584733
$ws.Cells.Item(13, 13).Value = 17.61728753847348
$ws.Cells.Item(13, 15).Value = 0.4970970511049586
$ws.Cells.Item(14, 3).Value = 0.04078051426184004
$ws.Cells.Item(14, 4).Value = 0.03791451430114989
$ws.Cells.Item(14, 5).Value = 4.258117322170165
$ws.Cells.Item(14, 6).Value = 0.2102194270677771
$ws.Cells.Item(14, 7).Value = 0.09886339014264678
$ws.Cells.Item(14, 8).Value = 0.2067535339812068
$ws.Cells.Item(14, 9).Value = 0.143268732803449
$ws.Cells.Item(14, 13).Value = 17.31531586811218
$ws.Cells.Item(14, 15).Value = 0.4979931985994028
$ws.Cells.Item(15, 3).Value = 0.04036926577546751
$ws.Cells.Item(15, 4).Value = 0.03751903239489707
$ws.Cells.Item(15, 5).Value = 4.211421892860585
$ws.Cells.Item(15, 6).Value = 0.2099027281106487
$ws.Cells.Item(15, 7).Value = 0.09844945656740123
$ws.Cells.Item(15, 8).Value = 0.2074847428234463
$ws.Cells.Item(15, 9).Value = 0.1431231490139595
$ws.Cells.Item(15, 13).Value = 17.13040587563086
$ws.Cells.Item(15, 15).Value = 0.4985945132945062
$ws.Cells.Item(16, 3).Value = 0.03800928050630148
$ws.Cells.Item(16, 4).Value = 0.03524982244211117
$ws.Cells.Item(16, 5).Value = 3.944400838201318
$ws.Cells.Item(16, 6).Value = 0.2083195150519188
$ws.Cells.Item(16, 7).Value = 0.0962944261520704
$ws.Cells.Item(16, 8).Value = 0.2118340729979238
$ws.Cells.Item(16, 9).Value = 0.1424642187642959
$ws.Cells.Item(16, 13).Value = 16.0715190177819
$ws.Cells.Item(16, 15).Value = 0.5028300983642993
$ws.Cells.Item(17, 3).Value = 0.03655863079694655
$ws.Cells.Item(17, 4).Value = 0.03385517139284389
$ws.Cells.Item(17, 5).Value = 3.781063447752217
$ws.Cells.Item(17, 6).Value = 0.2075509875231489
$ws.Cells.Item(17, 7).Value = 0.0951612295483244
$ws.Cells.Item(17, 8).Value = 0.2146421988306813
$ws.Cells.Item(17, 9).Value = 0.142213359492807
$ws.Cells.Item(17, 13).Value = 15.42249011217103
$ws.Cells.Item(17, 15).Value = 0.5061187200160759
$ws.Cells.Item(18, 3).Value = 0.03572322221403113
$ws.Cells.Item(18, 4).Value = 0.03305206039452457
$ws.Cells.Item(18, 5).Value = 3.68727450800975
$ws.Cells.Item(18, 6).Value = 0.2071818858973415
$ws.Cells.Item(18, 7).Value = 0.09457695076803674
$ws.Cells.Item(18, 8).Value = 0.2163080115822495
$ws.Cells.Item(18, 9).Value = 0.1421242227293718
$ws.Cells.Item(18, 13).Value = 15.04935014427383
$ws.Cells.Item(18, 15).Value = 0.5082582524339898
$ws.Cells.Item(19, 3).Value = 0.03544019203890514
$ws.Cells.Item(19, 4).Value = 0.03277997849943404
$ws.Cells.Item(19, 5).Value = 3.655545030664939
$ws.Cells.Item(19, 6).Value = 0.2070692943118004
$ws.Cells.Item(19, 7).Value = 0.09439054075308206
$ws.Cells.Item(19, 8).Value = 0.2168806537281043
$ws.Cells.Item(19, 9).Value = 0.1421034009708251
$ws.Cells.Item(19, 13).Value = 14.92303564542783
$ws.Cells.Item(19, 15).Value = 0.5090247264080432
$ws.Cells.Item(20, 3).Value = 0.03671316240904332
$ws.Cells.Item(20, 4).Value = 0.03400373211427166
$ws.Cells.Item(20, 5).Value = 3.798434271034694
$ws.Cells.Item(20, 6).Value = 0.2076252150937066
$ws.Cells.Item(20, 7).Value = 0.09527483133349079
$ws.Cells.Item(20, 8).Value = 0.2143380092122626
$ws.Cells.Item(20, 9).Value = 0.142234328857576
$ws.Cells.Item(20, 13).Value = 15.49156239960314
$ws.Cells.Item(20, 15).Value = 0.5057428467871858
$ws.Cells.Item(21, 3).Value = 0.04097766016256799
$ws.Cells.Item(21, 4).Value = 0.03810410815745513
$ws.Cells.Item(21, 5).Value = 4.280519719514018
$ws.Cells.Item(21, 6).Value = 0.2103754474755277
$ws.Cells.Item(21, 7).Value = 0.09906579091863676
$ws.Cells.Item(21, 8).Value = 0.2064057414837279
$ws.Cells.Item(21, 9).Value = 0.1433416900415878
$ws.Cells.Item(21, 13).Value = 17.40400025907445
$ws.Cells.Item(21, 15).Value = 0.4977190261681699
$ws.Cells.Item(22, 3).Value = 0.04375507095694786
$ws.Cells.Item(22, 4).Value = 0.04077560755401066
$ws.Cells.Item(22, 5).Value = 4.597325968700716
$ws.Cells.Item(22, 6).Value = 0.2128559721565253
$ws.Cells.Item(22, 7).Value = 0.1021852579053757
$ws.Cells.Item(22, 8).Value = 0.2016893295556201
$ws.Cells.Item(22, 9).Value = 0.1445825916842693
$ws.Cells.Item(22, 13).Value = 18.65629434069143
$ws.Cells.Item(22, 15).Value = 0.4948049197846274
$ws.Cells.Item(23, 3).Value = 0.04227362312272476
$ws.Cells.Item(23, 4).Value = 0.03935053424144996
$ws.Cells.Item(23, 5).Value = 4.428065137879116
$ws.Cells.Item(23, 6).Value = 0.211467908539376
$ws.Cells.Item(23, 7).Value = 0.1004595971672302
$ws.Cells.Item(23, 8).Value = 0.2041629408861922
$ws.Cells.Item(23, 9).Value = 0.1438716910348674
$ws.Cells.Item(23, 13).Value = 17.98765317611367
$ws.Cells.Item(23, 15).Value = 0.4961410460440163
$ws.Cells.Item(24, 3).Value = 0.03664330305595342
$ws.Cells.Item(24, 4).Value = 0.03393657188512123
$ws.Cells.Item(24, 5).Value = 3.790580564312052
$ws.Cells.Item(24, 6).Value = 0.2075914306242979
$ws.Cells.Item(24, 7).Value = 0.09522326306558426
$ws.Cells.Item(24, 8).Value = 0.2144753734731353
$ws.Cells.Item(24, 9).Value = 0.1422246773143954
$ws.Cells.Item(24, 13).Value = 15.4603348429095
$ws.Cells.Item(24, 15).Value = 0.505912005262843
$ws.Cells.Item(25, 3).Value = 0.03053090453801133
$ws.Cells.Item(25, 4).Value = 0.02806067395252398
$ws.Cells.Item(25, 5).Value = 3.108849405061079
$ws.Cells.Item(25, 6).Value = 0.2061798213195871
$ws.Cells.Item(25, 7).Value = 0.09213349231960422
$ws.Cells.Item(25, 8).Value = 0.2275262822814383
$ws.Cells.Item(25, 9).Value = 0.1425421912445408
$ws.Cells.Item(25, 13).Value = 12.74004103297898
$ws.Cells.Item(25, 15).Value = 0.5258662166592103
